$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced back to
# text (matching the workbook's original inlineStr cell type) - otherwise
# Excel silently reinterprets "206.56" etc. as a numeric value. Flagging the
# cell as Text before the write keeps the stored type a string; resetting the
# style to "Normal" afterwards drops the number-format override so the cell
# style index is left unchanged (matches the source which carries no explicit
# style on these cells).

$ws.Range('D2').Value = '27.168.56'
$ws.Range('E2').Value = '  -1.76%  '
$ws.Range('D3').Value = '1.559.40'
$ws.Range('E3').Value = '  -1.81%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '206.56'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('E6').Value = '  -1.13%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '22.23'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E10').Value = '  +0.22%  '
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('D12').Value = '1.780.87'
$ws.Range('E12').Value = '  -1.80%  '
$ws.Range('D13').Value = '1.554.22'
$ws.Range('E13').Value = '  -2.07%  '
$ws.Range('E14').Value = '  -2.09%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.516'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.84%  '
$ws.Range('E16').Value = '  -0.98%  '
$ws.Range('D17').Value = '27.149.26'
$ws.Range('E17').Value = '  -1.81%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '213.83'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.63%  '
$ws.Range('E19').Value = '  -1.31%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.23'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.21%  '
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('E22').Value = '  -0.47%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.37'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.24%  '
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '152.10'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.99%  '
$ws.Range('E26').Value = '  -3.19%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '14.89'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.51%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('E29').Value = '  -1.14%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.14'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.13%  '
$ws.Range('E31').Value = '  -1.34%  '
$ws.Range('E32').Value = '  -2.09%  '
$ws.Range('D33').Value = '1.383.35'
$ws.Range('E33').Value = '  +0.64%  '
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('E36').Value = '  -2.74%  '
$ws.Range('E37').Value = '  -1.37%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.813'
$ws.Range('D39').Style = "Normal"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.516'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.64%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.985'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.38%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.79'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +4.17%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.16'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '63.21'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.55%  '
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('D47').Value = '1.693.13'
$ws.Range('E47').Value = '  -1.79%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '85.49'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.33%  '
$ws.Range('D49').Value = '0.0₇0983'
$ws.Range('E49').Value = '  -2.41%  '
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('E51').Value = '  +0.14%  '
